$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "unit_number" column (F) moves from 0 to 1 for every data row (2:36).
$ws.Range("F2:F36").Value = 1

# The dropdown list restricting column A ("alk_lab" / water_type) to
# "nwfsc_oa, gagnon" is no longer needed and was removed.
$ws.Range("A2:A36").Validation.Delete()

# The user's cursor ended up parked on E2 when the file was saved.
$ws.Range("E2").Select()
